$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45190 -> 2023-09-21)
# that was updated to 45192 (-> 2023-09-23) for every data row (rows 2-535).
$ws.Range("C2:C535").Value = 45192
